$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- prime new shared strings in the exact order they were authored ---
$ws.Range("AH30").Value = "Tw"
$ws.Range("AA32").Value = "Te"
$ws.Range("W32").Value = "Li"
$ws.Range("AI27").Value = "Má"
$ws.Range("AJ27").Value = "Je"
$ws.Range("X30").Value = "Obs"
$ws.Range("Z30").Value = "pS"
$ws.Range("Y31").Value = "Nu"
$ws.Range("AI42").Value = "Tall"
$ws.Range("AI43").Value = "Wide"

# --- full pass: correct style + final value/formula per cell ---
$ws.Range("F2").Copy()
$ws.Range("R26").PasteSpecial(-4122)
$ws.Range("R26").Value = "E"
$ws.Range("F2").Copy()
$ws.Range("R27").PasteSpecial(-4122)
$ws.Range("R27").Value = 4
$ws.Range("J6").Copy()
$ws.Range("T27").PasteSpecial(-4122)
$ws.Range("T27").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("U27").PasteSpecial(-4122)
$ws.Range("U27").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("V27").PasteSpecial(-4122)
$ws.Range("V27").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("W27").PasteSpecial(-4122)
$ws.Range("W27").Value = "R"
$ws.Range("F2").Copy()
$ws.Range("X27").PasteSpecial(-4122)
$ws.Range("X27").Value = "LA"
$ws.Range("F2").Copy()
$ws.Range("Y27").PasteSpecial(-4122)
$ws.Range("Y27").Value = "LA"
$ws.Range("F2").Copy()
$ws.Range("Z27").PasteSpecial(-4122)
$ws.Range("Z27").Value = "LA"
$ws.Range("J6").Copy()
$ws.Range("AA27").PasteSpecial(-4122)
$ws.Range("AA27").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AB27").PasteSpecial(-4122)
$ws.Range("AB27").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("AC27").PasteSpecial(-4122)
$ws.Range("AC27").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("AD27").PasteSpecial(-4122)
$ws.Range("AD27").Value = "R"
$ws.Range("AI27").Value = "Má"
$ws.Range("AJ27").Value = "Je"
$ws.Range("R28").Value = 4
$ws.Range("L6").Copy()
$ws.Range("T28").PasteSpecial(-4122)
$ws.Range("T28").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("U28").PasteSpecial(-4122)
$ws.Range("U28").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("V28").PasteSpecial(-4122)
$ws.Range("V28").Value = "E"
$ws.Range("F2").Copy()
$ws.Range("W28").PasteSpecial(-4122)
$ws.Range("W28").Value = "LA"
$ws.Range("J6").Copy()
$ws.Range("X28").PasteSpecial(-4122)
$ws.Range("X28").Value = "R"
$ws.Range("F2").Copy()
$ws.Range("Y28").PasteSpecial(-4122)
$ws.Range("Y28").Value = "LA"
$ws.Range("J6").Copy()
$ws.Range("Z28").PasteSpecial(-4122)
$ws.Range("Z28").Value = "R"
$ws.Range("F2").Copy()
$ws.Range("AA28").PasteSpecial(-4122)
$ws.Range("AA28").Value = "LA"
$ws.Range("L6").Copy()
$ws.Range("AB28").PasteSpecial(-4122)
$ws.Range("AB28").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("AC28").PasteSpecial(-4122)
$ws.Range("AC28").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AD28").PasteSpecial(-4122)
$ws.Range("AD28").Value = "E"
$ws.Range("AH28").Value = "LA"
$ws.Range("AI28").Value = 10
$ws.Range("AJ28").Value = 10
$ws.Range("R29").Value = 4
$ws.Range("L6").Copy()
$ws.Range("T29").PasteSpecial(-4122)
$ws.Range("T29").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("U29").PasteSpecial(-4122)
$ws.Range("U29").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("V29").PasteSpecial(-4122)
$ws.Range("V29").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("W29").PasteSpecial(-4122)
$ws.Range("W29").Value = "R"
$ws.Range("X29").Value = "LA"
$ws.Range("Y29").Value = "LA"
$ws.Range("Z29").Value = "LA"
$ws.Range("J6").Copy()
$ws.Range("AA29").PasteSpecial(-4122)
$ws.Range("AA29").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("AB29").PasteSpecial(-4122)
$ws.Range("AB29").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AC29").PasteSpecial(-4122)
$ws.Range("AC29").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("AD29").PasteSpecial(-4122)
$ws.Range("AD29").Value = "E"
$ws.Range("AH29").Value = "E"
$ws.Range("AI29").Value = 60
$ws.Range("AJ29").Formula = "=SUM(R27:R39)"
$ws.Range("R30").Value = 6
$ws.Range("L6").Copy()
$ws.Range("T30").PasteSpecial(-4122)
$ws.Range("T30").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("U30").PasteSpecial(-4122)
$ws.Range("U30").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("V30").PasteSpecial(-4122)
$ws.Range("V30").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("W30").PasteSpecial(-4122)
$ws.Range("W30").Value = "R"
$ws.Range("AO7").Copy()
$ws.Range("X30").PasteSpecial(-4122)
$ws.Range("X30").Value = "Obs"
$ws.Range("F2").Copy()
$ws.Range("Y30").PasteSpecial(-4122)
$ws.Range("Y30").Value = "LA"
$ws.Range("AO7").Copy()
$ws.Range("Z30").PasteSpecial(-4122)
$ws.Range("Z30").Value = "pS"
$ws.Range("J6").Copy()
$ws.Range("AA30").PasteSpecial(-4122)
$ws.Range("AA30").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AB30").PasteSpecial(-4122)
$ws.Range("AB30").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("AC30").PasteSpecial(-4122)
$ws.Range("AC30").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("AD30").PasteSpecial(-4122)
$ws.Range("AD30").Value = "E"
$ws.Range("AH30").Value = "Tw"
$ws.Range("AI30").Value = 6
$ws.Range("AJ30").Value = 6
$ws.Range("R31").Value = 6
$ws.Range("J6").Copy()
$ws.Range("T31").PasteSpecial(-4122)
$ws.Range("T31").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("U31").PasteSpecial(-4122)
$ws.Range("U31").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("V31").PasteSpecial(-4122)
$ws.Range("V31").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("W31").PasteSpecial(-4122)
$ws.Range("W31").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("X31").PasteSpecial(-4122)
$ws.Range("X31").Value = "R"
$ws.Range("AO7").Copy()
$ws.Range("Y31").PasteSpecial(-4122)
$ws.Range("Y31").Value = "Nu"
$ws.Range("J6").Copy()
$ws.Range("Z31").PasteSpecial(-4122)
$ws.Range("Z31").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AA31").PasteSpecial(-4122)
$ws.Range("AA31").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("AB31").PasteSpecial(-4122)
$ws.Range("AB31").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("AC31").PasteSpecial(-4122)
$ws.Range("AC31").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("AD31").PasteSpecial(-4122)
$ws.Range("AD31").Value = "R"
$ws.Range("R32").Value = 4
$ws.Range("L6").Copy()
$ws.Range("T32").PasteSpecial(-4122)
$ws.Range("T32").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("U32").PasteSpecial(-4122)
$ws.Range("U32").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("V32").PasteSpecial(-4122)
$ws.Range("V32").Value = "E"
$ws.Range("W32").Value = "Li"
$ws.Range("J6").Copy()
$ws.Range("X32").PasteSpecial(-4122)
$ws.Range("X32").Value = "R"
$ws.Range("Y32").Value = "Tw"
$ws.Range("J6").Copy()
$ws.Range("Z32").PasteSpecial(-4122)
$ws.Range("Z32").Value = "R"
$ws.Range("AA32").Value = "Te"
$ws.Range("L6").Copy()
$ws.Range("AB32").PasteSpecial(-4122)
$ws.Range("AB32").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("AC32").PasteSpecial(-4122)
$ws.Range("AC32").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AD32").PasteSpecial(-4122)
$ws.Range("AD32").Value = "E"
$ws.Range("R33").Value = 2
$ws.Range("J6").Copy()
$ws.Range("T33").PasteSpecial(-4122)
$ws.Range("T33").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("U33").PasteSpecial(-4122)
$ws.Range("U33").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("V33").PasteSpecial(-4122)
$ws.Range("V33").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("W33").PasteSpecial(-4122)
$ws.Range("W33").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("X33").PasteSpecial(-4122)
$ws.Range("X33").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("Y33").PasteSpecial(-4122)
$ws.Range("Y33").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("Z33").PasteSpecial(-4122)
$ws.Range("Z33").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("AA33").PasteSpecial(-4122)
$ws.Range("AA33").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("AB33").PasteSpecial(-4122)
$ws.Range("AB33").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AC33").PasteSpecial(-4122)
$ws.Range("AC33").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("AD33").PasteSpecial(-4122)
$ws.Range("AD33").Value = "R"
$ws.Range("R34").Value = 4
$ws.Range("L6").Copy()
$ws.Range("T34").PasteSpecial(-4122)
$ws.Range("T34").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("U34").PasteSpecial(-4122)
$ws.Range("U34").Value = "E"
$ws.Range("V34").Value = "Tw"
$ws.Range("W34").Value = "Tw"
$ws.Range("J6").Copy()
$ws.Range("X34").PasteSpecial(-4122)
$ws.Range("X34").Value = "R"
$ws.Range("AO7").Copy()
$ws.Range("Y34").PasteSpecial(-4122)
$ws.Range("Y34").Value = "ST"
$ws.Range("J6").Copy()
$ws.Range("Z34").PasteSpecial(-4122)
$ws.Range("Z34").Value = "R"
$ws.Range("AA34").Value = "Tw"
$ws.Range("AB34").Value = "Tw"
$ws.Range("L6").Copy()
$ws.Range("AC34").PasteSpecial(-4122)
$ws.Range("AC34").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("AD34").PasteSpecial(-4122)
$ws.Range("AD34").Value = "E"
$ws.Range("R35").Value = 2
$ws.Range("J6").Copy()
$ws.Range("T35").PasteSpecial(-4122)
$ws.Range("T35").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("U35").PasteSpecial(-4122)
$ws.Range("U35").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("V35").PasteSpecial(-4122)
$ws.Range("V35").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("W35").PasteSpecial(-4122)
$ws.Range("W35").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("X35").PasteSpecial(-4122)
$ws.Range("X35").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("Y35").PasteSpecial(-4122)
$ws.Range("Y35").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("Z35").PasteSpecial(-4122)
$ws.Range("Z35").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("AA35").PasteSpecial(-4122)
$ws.Range("AA35").Value = "R"
$ws.Range("J6").Copy()
$ws.Range("AB35").PasteSpecial(-4122)
$ws.Range("AB35").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AC35").PasteSpecial(-4122)
$ws.Range("AC35").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("AD35").PasteSpecial(-4122)
$ws.Range("AD35").Value = "R"
$ws.Range("R36").Value = 4
$ws.Range("L6").Copy()
$ws.Range("T36").PasteSpecial(-4122)
$ws.Range("T36").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("U36").PasteSpecial(-4122)
$ws.Range("U36").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("V36").PasteSpecial(-4122)
$ws.Range("V36").Value = "E"
$ws.Range("W36").Value = "S2"
$ws.Range("J6").Copy()
$ws.Range("X36").PasteSpecial(-4122)
$ws.Range("X36").Value = "R"
$ws.Range("Y36").Value = "S1"
$ws.Range("J6").Copy()
$ws.Range("Z36").PasteSpecial(-4122)
$ws.Range("Z36").Value = "R"
$ws.Range("AA36").Value = "S3"
$ws.Range("L6").Copy()
$ws.Range("AB36").PasteSpecial(-4122)
$ws.Range("AB36").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("AC36").PasteSpecial(-4122)
$ws.Range("AC36").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AD36").PasteSpecial(-4122)
$ws.Range("AD36").Value = "E"
$ws.Range("R37").Value = 6
$ws.Range("J6").Copy()
$ws.Range("T37").PasteSpecial(-4122)
$ws.Range("T37").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("U37").PasteSpecial(-4122)
$ws.Range("U37").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("V37").PasteSpecial(-4122)
$ws.Range("V37").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("W37").PasteSpecial(-4122)
$ws.Range("W37").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("X37").PasteSpecial(-4122)
$ws.Range("X37").Value = "R"
$ws.Range("Y37").Value = "Tw"
$ws.Range("J6").Copy()
$ws.Range("Z37").PasteSpecial(-4122)
$ws.Range("Z37").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AA37").PasteSpecial(-4122)
$ws.Range("AA37").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("AB37").PasteSpecial(-4122)
$ws.Range("AB37").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("AC37").PasteSpecial(-4122)
$ws.Range("AC37").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("AD37").PasteSpecial(-4122)
$ws.Range("AD37").Value = "R"
$ws.Range("R38").Value = 9
$ws.Range("L6").Copy()
$ws.Range("T38").PasteSpecial(-4122)
$ws.Range("T38").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("U38").PasteSpecial(-4122)
$ws.Range("U38").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("V38").PasteSpecial(-4122)
$ws.Range("V38").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("W38").PasteSpecial(-4122)
$ws.Range("W38").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("X38").PasteSpecial(-4122)
$ws.Range("X38").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("Y38").PasteSpecial(-4122)
$ws.Range("Y38").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("Z38").PasteSpecial(-4122)
$ws.Range("Z38").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("AA38").PasteSpecial(-4122)
$ws.Range("AA38").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AB38").PasteSpecial(-4122)
$ws.Range("AB38").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("AC38").PasteSpecial(-4122)
$ws.Range("AC38").Value = "E"
$ws.Range("L6").Copy()
$ws.Range("AD38").PasteSpecial(-4122)
$ws.Range("AD38").Value = "E"
$ws.Range("R39").Value = 5
$ws.Range("L6").Copy()
$ws.Range("U39").PasteSpecial(-4122)
$ws.Range("U39").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("V39").PasteSpecial(-4122)
$ws.Range("V39").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("W39").PasteSpecial(-4122)
$ws.Range("W39").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("X39").PasteSpecial(-4122)
$ws.Range("X39").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("Y39").PasteSpecial(-4122)
$ws.Range("Y39").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("Z39").PasteSpecial(-4122)
$ws.Range("Z39").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AA39").PasteSpecial(-4122)
$ws.Range("AA39").Value = "E"
$ws.Range("J6").Copy()
$ws.Range("AB39").PasteSpecial(-4122)
$ws.Range("AB39").Value = "R"
$ws.Range("L6").Copy()
$ws.Range("AC39").PasteSpecial(-4122)
$ws.Range("AC39").Value = "E"
$ws.Range("AH42").Value = 13
$ws.Range("AI42").Value = "Tall"
$ws.Range("AH43").Value = 11
$ws.Range("AI43").Value = "Wide"

$excel.CutCopyMode = 0

# --- view state ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("AF27").Select()
